# Scheduled-runner update: refresh computed market/profit columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) across the leve sheets
# with newly scraped marketboard data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 1499.6666
$ws.Range("J21").Value = 1499.6666
$ws.Range("L21").Value = 1499.6666
$ws.Range("N21").Value = -2435.6666

$ws.Range("H23").Value = 1499.6666
$ws.Range("J23").Value = 1499.6666
$ws.Range("L23").Value = 1499.6666
$ws.Range("N23").Value = -1967.6666

$ws.Range("H70").Value = 2304.0588
$ws.Range("I70").Value = 1926.6
$ws.Range("J70").Value = 2843.2856
$ws.Range("K70").Value = 5779.799999999999
$ws.Range("L70").Value = 8529.856800000001
$ws.Range("M70").Value = -5509.799999999999
$ws.Range("N70").Value = -9069.856800000001

$ws.Range("H73").Value = 2304.0588
$ws.Range("I73").Value = 1926.6
$ws.Range("J73").Value = 2843.2856
$ws.Range("K73").Value = 5779.799999999999
$ws.Range("L73").Value = 8529.856800000001
$ws.Range("M73").Value = -4843.799999999999
$ws.Range("N73").Value = -10401.8568

$ws.Range("H113").Value = 4844.3
$ws.Range("I113").Value = 3775.923
$ws.Range("K113").Value = 3775.923
$ws.Range("M113").Value = -521.9229999999998

$ws.Range("H132").Value = 3575.8064
$ws.Range("I132").Value = 3347.4138
$ws.Range("K132").Value = 10042.2414
$ws.Range("M132").Value = -7512.241399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 29423452
$ws.Range("J32").Value = 20948.084
$ws.Range("L32").Value = 20948.084
$ws.Range("N32").Value = -21522.084

$ws.Range("H61").Value = 83338800
$ws.Range("J61").Value = 9252.333000000001
$ws.Range("L61").Value = 9252.333000000001
$ws.Range("N61").Value = -9676.333000000001

$ws.Range("H74").Value = 40046944
$ws.Range("I74").Value = 58890750
$ws.Range("J74").Value = 3851.75
$ws.Range("K74").Value = 58890750
$ws.Range("L74").Value = 3851.75
$ws.Range("M74").Value = -58889876
$ws.Range("N74").Value = -5599.75

$ws.Range("H77").Value = 40046944
$ws.Range("I77").Value = 58890750
$ws.Range("J77").Value = 3851.75
$ws.Range("K77").Value = 294453750
$ws.Range("L77").Value = 19258.75
$ws.Range("M77").Value = -294449382
$ws.Range("N77").Value = -27994.75

$ws.Range("H102").Value = 2375.6
$ws.Range("I102").Value = 1468.125
$ws.Range("K102").Value = 1468.125
$ws.Range("M102").Value = 153.875

$ws.Range("H110").Value = 12155
$ws.Range("I110").Value = 14285.296
$ws.Range("K110").Value = 14285.296
$ws.Range("M110").Value = -12240.296

$ws.Range("H113").Value = 74497.5
$ws.Range("J113").Value = 74497.5
$ws.Range("L113").Value = 74497.5
$ws.Range("N113").Value = -83175.5

$ws.Range("H136").Value = 83338800
$ws.Range("J136").Value = 9252.333000000001
$ws.Range("L136").Value = 27756.999
$ws.Range("N136").Value = -32856.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 3127.7856
$ws.Range("I99").Value = 1753.2222
$ws.Range("K99").Value = 1753.2222
$ws.Range("M99").Value = -255.2221999999999

$ws.Range("H105").Value = 1210
$ws.Range("I105").Value = 1137.5
$ws.Range("K105").Value = 1137.5
$ws.Range("M105").Value = 609.5

$ws.Range("H107").Value = 3015.2927
$ws.Range("J107").Value = 5357.9165
$ws.Range("L107").Value = 5357.9165
$ws.Range("N107").Value = -9197.916499999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4021899.8
$ws.Range("J4").Value = 22750
$ws.Range("L4").Value = 22750
$ws.Range("N4").Value = -22974

$ws.Range("H31").Value = 22227332
$ws.Range("I31").Value = 3918.2593
$ws.Range("J31").Value = 55562456
$ws.Range("K31").Value = 3918.2593
$ws.Range("L31").Value = 55562456
$ws.Range("M31").Value = -3623.2593
$ws.Range("N31").Value = -55563046

$ws.Range("H34").Value = 22227332
$ws.Range("I34").Value = 3918.2593
$ws.Range("J34").Value = 55562456
$ws.Range("K34").Value = 3918.2593
$ws.Range("L34").Value = 55562456
$ws.Range("M34").Value = -3716.2593
$ws.Range("N34").Value = -55562860

$ws.Range("H99").Value = 6035.7085
$ws.Range("I99").Value = 7139.5293
$ws.Range("K99").Value = 7139.5293
$ws.Range("M99").Value = -5641.5293

$ws.Range("H105").Value = 10780
$ws.Range("I105").Value = 1459.5
$ws.Range("J105").Value = 21964.6
$ws.Range("K105").Value = 1459.5
$ws.Range("L105").Value = 21964.6
$ws.Range("M105").Value = 287.5
$ws.Range("N105").Value = -25458.6

$ws.Range("H107").Value = 1222.2307
$ws.Range("J107").Value = 1417
$ws.Range("L107").Value = 1417
$ws.Range("N107").Value = -5257

$ws.Range("H126").Value = 6035.7085
$ws.Range("I126").Value = 7139.5293
$ws.Range("K126").Value = 21418.5879
$ws.Range("M126").Value = -18948.5879

$ws.Range("H132").Value = 3909.0356
$ws.Range("I132").Value = 3855.6086
$ws.Range("J132").Value = 4154.8
$ws.Range("K132").Value = 11566.8258
$ws.Range("L132").Value = 12464.4
$ws.Range("M132").Value = -9036.825800000001
$ws.Range("N132").Value = -17524.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H61").Value = 374.8
$ws.Range("I61").Value = 423.75
$ws.Range("K61").Value = 1271.25
$ws.Range("M61").Value = -1056.25

$ws.Range("H103").Value = 243.66667
$ws.Range("I103").Value = 243.66667
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 731.00001
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 147.99999
$ws.Range("N103").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 3463.3333
$ws.Range("I97").Value = 1099.5
$ws.Range("J97").Value = 4138.7144
$ws.Range("K97").Value = 1099.5
$ws.Range("L97").Value = 4138.7144
$ws.Range("M97").Value = -603.5
$ws.Range("N97").Value = -5130.7144

$ws.Range("H102").Value = 4221.4443
$ws.Range("I102").Value = 2994.5
$ws.Range("K102").Value = 2994.5
$ws.Range("M102").Value = -1372.5

$ws.Range("H107").Value = 441.52942
$ws.Range("I107").Value = 589.6667
$ws.Range("J107").Value = 274.875
$ws.Range("K107").Value = 589.6667
$ws.Range("L107").Value = 274.875
$ws.Range("M107").Value = 1330.3333
$ws.Range("N107").Value = -4114.875

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2548
$ws.Range("I22").Value = 1813.6428
$ws.Range("J22").Value = 3062.05
$ws.Range("K22").Value = 1813.6428
$ws.Range("L22").Value = 3062.05
$ws.Range("M22").Value = -1518.6428
$ws.Range("N22").Value = -3652.05

$ws.Range("H27").Value = 2548
$ws.Range("I27").Value = 1813.6428
$ws.Range("J27").Value = 3062.05
$ws.Range("K27").Value = 1813.6428
$ws.Range("L27").Value = 3062.05
$ws.Range("M27").Value = -1706.6428
$ws.Range("N27").Value = -3276.05

$ws.Range("H82").Value = 3235.4
$ws.Range("I82").Value = 1977.6666
$ws.Range("K82").Value = 1977.6666
$ws.Range("M82").Value = -1616.6666

$ws.Range("H85").Value = 3235.4
$ws.Range("I85").Value = 1977.6666
$ws.Range("K85").Value = 1977.6666
$ws.Range("M85").Value = -729.6666

$ws.Range("H100").Value = 3099.5
$ws.Range("I100").Value = 1764.3334
$ws.Range("K100").Value = 1764.3334
$ws.Range("M100").Value = -1223.3334

$ws.Range("H132").Value = 125001950
$ws.Range("I132").Value = 1871.8572
$ws.Range("J132").Value = 222224240
$ws.Range("K132").Value = 5615.571599999999
$ws.Range("L132").Value = 666672720
$ws.Range("M132").Value = -3085.571599999999
$ws.Range("N132").Value = -666677780

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 622
$ws.Range("I107").Value = 473
$ws.Range("J107").Value = 882.75
$ws.Range("K107").Value = 1419
$ws.Range("L107").Value = 2648.25
$ws.Range("M107").Value = 501
$ws.Range("N107").Value = -6488.25
